$d = $word.ActiveDocument

# Helper: set the text of [rangeStart, rangeEnd) and then force it into
# its own <w:r> by toggling a formatting property on and back off again
# (the engine coalesces adjacent runs with identical rPr, so the text
# assignment must happen *before* the toggle or the split gets merged
# back out).
function Set-AndSplitRun($rangeStart, $rangeEnd, $text) {
    $sub = $d.Range($rangeStart, $rangeEnd)
    $sub.Text = $text
    $sub2 = $d.Range($rangeStart, $rangeEnd)
    $sub2.Bold = 1
    $sub2.Bold = 0
}

# ---------------------------------------------------------------------
# 1) Fighting Profiles table (first table): "Analyzer" row, Shoot column
#    "RW+2" -> "R" / "S" / "+2"
# ---------------------------------------------------------------------
$t1 = $d.Tables(1)
$cell1 = $t1.Cell(2, 4)
$start1 = $cell1.Range.Start
$a1 = $start1 + 1
$b1 = $start1 + 2
Set-AndSplitRun $a1 $b1 "S"

# ---------------------------------------------------------------------
# 2) Fighting Profiles table: "Predictor" row, Strike column
#    "MW+0" -> "M" / "S" / "+0"
# ---------------------------------------------------------------------
$cell2 = $t1.Cell(3, 2)
$start2 = $cell2.Range.Start
$a2 = $start2 + 1
$b2 = $start2 + 2
Set-AndSplitRun $a2 $b2 "S"

# ---------------------------------------------------------------------
# 3) Fighting Profiles table: "Predictor" row, Shoot column
#    "RW+3" -> "R" / "W" / (bookmark _GoBack) / "+3"
# ---------------------------------------------------------------------
$cell3 = $t1.Cell(3, 4)
$start3 = $cell3.Range.Start
$a3 = $start3 + 1
$b3 = $start3 + 2
Set-AndSplitRun $a3 $b3 "W"

$bmStart = $start3 + 2
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

$c3 = $start3 + 2
$e3 = $start3 + 4
Set-AndSplitRun $c3 $e3 "+3"

# ---------------------------------------------------------------------
# 4) Remove the old _GoBack bookmark sitting after the "None" run
#    (alt-init-system list item).
# ---------------------------------------------------------------------
$noneSearch = $d.Content
$noneSearch.Find.Execute("None", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$noneParaRange = $noneSearch.Paragraphs(1).Range
$cleanParaXml = '<?xml version="1.0"?>' + `
 '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
 '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
 '<pkg:xmlData>' + `
 '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
 '<w:body>' + `
 '<w:p>' + `
 '<w:pPr>' + `
 '<w:pStyle w:val="ListParagraph"/>' + `
 '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
 '<w:spacing w:line="276" w:lineRule="auto"/>' + `
 '<w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr>' + `
 '</w:pPr>' + `
 '<w:r>' + `
 '<w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr>' + `
 '<w:t>None</w:t>' + `
 '</w:r>' + `
 '</w:p>' + `
 '</w:body></w:document>' + `
 '</pkg:xmlData></pkg:part></pkg:package>'
$noneParaRange.InsertXML($cleanParaXml)

# ---------------------------------------------------------------------
# 5) "Investigation +2" -> "Investigation +" / "4"
# ---------------------------------------------------------------------
$invRange = $d.Content
$invRange.Find.Execute("Investigation +2") | Out-Null
$invEnd = $invRange.End
$a5 = $invEnd - 1
Set-AndSplitRun $a5 $invEnd "4"

# ---------------------------------------------------------------------
# 6) "Science/Math +2" -> "Science/Math +" / "4"
# ---------------------------------------------------------------------
$sciRange = $d.Content
$sciRange.Find.Execute("Science/Math +2") | Out-Null
$sciEnd = $sciRange.End
$a6 = $sciEnd - 1
Set-AndSplitRun $a6 $sciEnd "4"
